# Actualización 10 de Mayo
# Updates the "BRINDA SOPORTE TÉCNICO DE MANERA PRESENCIAL" (row 7) statistics
# for Sánchez Sánchez Miguel / 4ASV in each of the three partial-exam sheets.

$wb = $excel.ActiveWorkbook

# --- Hoja "1er Parcial" ---
$ws = $wb.Worksheets.Item("1er Parcial")
$ws.Range("E7").Value = 30
$ws.Range("F7").Value = 6
$ws.Range("G7").Value = 83.33
$ws.Range("H7").Value = 16.67
$ws.Range("I7").Value = 8.300000000000001
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0

# --- Hoja "2o Parcial" ---
$ws = $wb.Worksheets.Item("2o Parcial")
$ws.Range("E7").Value = 31
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 86.11
$ws.Range("H7").Value = 13.89
$ws.Range("I7").Value = 8
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0

# --- Hoja "3er Parcial" ---
$ws = $wb.Worksheets.Item("3er Parcial")
$ws.Range("E7").Value = 31
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 86.11
$ws.Range("H7").Value = 13.89
$ws.Range("I7").Value = 8.4
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
